# Review of code, small bugs fixed
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Sheet1: update ticker rows 2 & 3, remove rows 4-6 (amzn, celh, rivn) ---

# Row 2: aapl -> tsla
$ws1.Range("A2").Value = "tsla"
$ws1.Range("B2").Value = 12
$ws1.Range("C2").Value = 397.2099914550781
$ws1.Range("D2").Value = 417.0700073242188
$ws1.Range("E2").Value = 5.093026716730011

# Row 3: goog -> nmm
$ws1.Range("A3").Value = "nmm"
$ws1.Range("B3").Value = 37
$ws1.Range("C3").Value = 55.43999862670898
$ws1.Range("D3").Value = 60.33000183105469
$ws1.Range("E3").Value = 11.72438701480872

# Remove trailing rows 4 (amzn), 5 (celh), 6 (rivn) - delete from bottom up
$ws1.Rows.Item(6).Delete()
$ws1.Rows.Item(5).Delete()
$ws1.Rows.Item(4).Delete()

# --- Summary sheet: update aggregate row 2 ---
$ws2.Range("B2").Value = 2
$ws2.Range("C2").Value = 6817.79984664917
$ws2.Range("D2").Value = 49
$ws2.Range("E2").Value = 139.1387723805953
$ws2.Range("F2").Value = 5.46780506258587
